# Database Schema minor change
# Applies datatype/attribute corrections across the schema sheets:
# - Many ID / numeric columns change Datatype from Varchar/Number/number-int to int
# - Several Attribute cells gain "not null" / "unique" constraints
# - A few label corrections (House Number -> House number, Phone Number -> Phone number,
#   User ID -> Provider ID / Favourite Pro ID / Block Pro ID, Rate -> Rating)
# - The stray empty "Service Provider" row on the Book sheet is removed
# - Selections (active cell) on each sheet are updated

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "User"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("User")
$ws.Range("B3").Value = "int"
$ws.Range("C4").Value = "Not null,unique"
$ws.Range("B5").Value = "int"
$ws.Range("C5").Value = "not null"
$ws.Range("C6").Value = "not null"
$ws.Range("C7").Value = "not null"
$ws.Range("C10").Value = "not null"
$ws.Activate() | Out-Null
$ws.Range("B5").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "Service Provider"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Service Provider")
$ws.Range("B3").Value = "int"
$ws.Range("C4").Value = "not null"
$ws.Range("C5").Value = "not null"
$ws.Range("C6").Value = "not null,unique"
$ws.Range("B7").Value = "Int"
$ws.Range("C7").Value = "not null"
$ws.Range("C8").Value = "not null"
$ws.Activate() | Out-Null
$ws.Range("C9").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "Book"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Book")
# Remove the stray blank "Service Provider" row (old row 11)
$ws.Rows("11:11").Delete() | Out-Null
$ws.Range("B3").Value = "int"
$ws.Range("B4").Value = "int"
$ws.Range("A5").Value = "Provider ID"
$ws.Range("B5").Value = "int"
$ws.Range("C6").Value = "not null"
$ws.Range("C7").Value = "not null"
$ws.Range("A14").Value = "Rating"
$ws.Activate() | Out-Null
$ws.Range("C14").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "Favourite Pros"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Favourite Pros")
$ws.Range("B3").Value = "int"
$ws.Range("A4").Value = "Favourite Pro ID"
$ws.Range("B4").Value = "int"
$ws.Range("B5").Value = "int"

# ---------------------------------------------------------------
# Sheet "Address"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Address")
$ws.Range("B3").Value = "int"
$ws.Range("B4").Value = "int"
$ws.Range("C5").Value = "not null"
$ws.Range("A6").Value = "House number"
$ws.Range("B6").Value = "int"
$ws.Range("B7").Value = "int"
$ws.Range("C7").Value = "not null"
$ws.Range("C8").Value = "not null"
$ws.Range("A9").Value = "Phone number"
$ws.Range("B9").Value = "int"
$ws.Activate() | Out-Null
$ws.Range("C9").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "Block"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Block")
$ws.Range("B3").Value = "int"
$ws.Range("A4").Value = "Block Pro ID"
$ws.Range("B4").Value = "int"
$ws.Range("B5").Value = "int"
$ws.Activate() | Out-Null
$ws.Range("B7").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "Extra Services"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Extra Services")
$ws.Range("B3").Value = "int"
$ws.Range("C4").Value = "not null"
$ws.Range("C5").Value = "not null"
$ws.Range("C6").Value = "not null"
$ws.Activate() | Out-Null
$ws.Range("C5").Select() | Out-Null
